# 0191｜EASY｜Number of 1 Bits
# - Convert L73 from the stray text "8.2 " shared-string into a real number (8.2)
# - Append 3 new rows (74-76) documenting the 3 solution methods for
#   LeetCode 0191 "Number of 1 Bits" (Bit Manipulation)
# - Move the active selection down to the newly added data (F74)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix L73: was stored as text "8.2 ", should be the numeric value 8.2 ---
$ws.Range("L73").Value = 8.2

# --- Copy the formatting of the last existing row (73) down onto the three
#     new rows so fonts/fills/borders/number-formats match the rest of the
#     table exactly ---
$ws.Range("A73:P73").Copy()
$ws.Range("A74:P76").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 74: method1.cpp / bitset ---
$ws.Range("A74").Value = "0191"
$ws.Range("B74").Value = "EASY"
$ws.Range("C74").Value = "Number of 1 Bits"
$ws.Range("D74").Value = "method1.cpp"
$ws.Range("E74").Value = "Bit Manipulation"
$ws.Range("F74").Value = "bitset"
$ws.Range("G74").Value = "DONE"
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 100
$ws.Range("J74").Value = 4
$ws.Range("K74").Value = 66.24
$ws.Range("L74").Value = 8.3
$ws.Range("M74").Value = 60.98
$ws.Range("N74").Value = 43847
$ws.Range("O74").Value = 0.64930555555555558
$ws.Range("P74").Value = $false

# --- Row 75: method2.cpp / 换算2进制 ---
$ws.Range("A75").Value = "0191"
$ws.Range("B75").Value = "EASY"
$ws.Range("C75").Value = "Number of 1 Bits"
$ws.Range("D75").Value = "method2.cpp"
$ws.Range("E75").Value = "Bit Manipulation"
$ws.Range("F75").Value = "换算2进制"
$ws.Range("G75").Value = "DONE"
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 100
$ws.Range("J75").Value = 4
$ws.Range("K75").Value = 66.24
$ws.Range("L75").Value = 8.2
$ws.Range("M75").Value = 78.05
$ws.Range("N75").Value = 43847
$ws.Range("O75").Value = 0.64930555555555558
$ws.Range("P75").Value = $false

# --- Row 76: method3.cpp / bitwise operator ---
$ws.Range("A76").Value = "0191"
$ws.Range("B76").Value = "EASY"
$ws.Range("C76").Value = "Number of 1 Bits"
$ws.Range("D76").Value = "method3.cpp"
$ws.Range("E76").Value = "Bit Manipulation"
$ws.Range("F76").Value = "bitwise operator"
$ws.Range("G76").Value = "DONE"
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 100
$ws.Range("J76").Value = 4
$ws.Range("K76").Value = 66.24
$ws.Range("L76").Value = 8.3
$ws.Range("M76").Value = 60.98
$ws.Range("N76").Value = 43847
$ws.Range("O76").Value = 0.64930555555555558
$ws.Range("P76").Value = $false

# --- Move the selection to match where the author ended up editing ---
$ws.Range("F74").Select()
